# Append/refresh timestamp: 2025-09-29 18:32 JST
# Updates the "取得日時" (acquisition timestamp) column A for all data rows
# on the "ランサーズ" sheet from "2025-09-29 18:26:32" to "2025-09-29 18:32:37".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-29 18:32:37"

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value -ne $null -and $cell.Value -ne "") {
        $cell.Value = $newTimestamp
    }
}
